# Logged Week 17 data and fixed Simulate_Season.py tiebreaking method
# This script appends the Week 17 per-play yardage / special-teams logs
# (shared strings) and updates the season-to-date aggregate totals on the
# YDS, OFF, DEF, ST, TURNS and PEN sheets.

$wb = $excel.ActiveWorkbook

function Append-Values {
    param(
        $Sheet,
        [string]$CellRef,
        [string]$ToAppend
    )
    $range = $Sheet.Range($CellRef)
    $current = [string]$range.Value2
    $range.Value = $current + $ToAppend
}

# ---------------------------------------------------------------------
# Sheet "YDS" - append Week 17 per-play yardage logs
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("YDS")
Append-Values $ws "B2" " 4 1 4 15 7 2 6 9 2 9 7 0 7 5 4 8 2 4 4 7 10 8 -1 4 1 13 5 3 7 4 6 4 8"
Append-Values $ws "B3" " 8 7 12 6 14 9 17 8 6 8 8 4 12 20 20 17 8 8 8 7 7 13 46 8 9"
Append-Values $ws "C2" " 4 15 5 2 5 1 5 1 4 2 4 3 4 9 5 4 7"
Append-Values $ws "C3" " 7 9 40 2 8 24 11 22 5 4 3 -5 4 16 6 9 28"

# ---------------------------------------------------------------------
# Sheet "OFF" - season totals through Week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 209
$ws.Range("D2").Value = 13
$ws.Range("F2").Value = 54
$ws.Range("G2").Value = 64
$ws.Range("J2").Value = 37
$ws.Range("L2").Value = 318
$ws.Range("M2").Value = 220
$ws.Range("Q2").Value = 561

$ws.Range("C3").Value = 149
$ws.Range("D3").Value = 5
$ws.Range("E3").Value = 34
$ws.Range("F3").Value = 93
$ws.Range("G3").Value = 23
$ws.Range("I3").Value = 52
$ws.Range("J3").Value = 53

# ---------------------------------------------------------------------
# Sheet "DEF" - season totals through Week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 195
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 59
$ws.Range("G2").Value = 59
$ws.Range("I2").Value = 6
$ws.Range("L2").Value = 314
$ws.Range("M2").Value = 191
$ws.Range("Q2").Value = 568

$ws.Range("C3").Value = 155
$ws.Range("E3").Value = 28
$ws.Range("F3").Value = 102
$ws.Range("G3").Value = 37
$ws.Range("H3").Value = 29
$ws.Range("J3").Value = 49
$ws.Range("N3").Value = 17

# ---------------------------------------------------------------------
# Sheet "ST" - append Week 17 logs and update season totals
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 98
$ws.Range("D2").Value = 47
$ws.Range("F2").Value = 435
$ws.Range("G2").Value = 423
$ws.Range("H2").Value = 4
$ws.Range("L2").Value = 120
$ws.Range("M2").Value = 89

$ws.Range("B3").Value = 41

Append-Values $ws "B4" " 67 61 62 60 61"
Append-Values $ws "D3" " 40"
Append-Values $ws "D4" " 3"
Append-Values $ws "B5" " 25 18 16 9 20"
Append-Values $ws "D5" " 0 0 0 0"
Append-Values $ws "B6" " 14 27 9"

# ---------------------------------------------------------------------
# Sheet "TURNS" - season totals through Week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("C2").Value = 9
$ws.Range("D3").Value = 7

# ---------------------------------------------------------------------
# Sheet "PEN" - season totals through Week 17
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B3").Value = 21
